$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the content of row 16 and row 17 (two species observation
# records traded places), while a handful of columns that happened to share
# identical values in both rows (D, I, T, U, V, W, Y, AA, AD, AE, AG, AT, AY)
# remain unchanged either way.

# --- Row 16 becomes what Row 17 used to be ---
$ws.Range("A16").Value2 = 131066881
$ws.Range("B16").Value2 = 57884
$ws.Range("E16").Value2 = 100109
$ws.Range("F16").Value2 = "Tretåig hackspett"
$ws.Range("G16").Value2 = "Picoides tridactylus"
$ws.Range("H16").Value2 = "(Linnaeus, 1758)"
$ws.Range("M16").Value2 = "äldre spår"
$ws.Range("P16").Value2 = "Färntjärnen, Vrm"
$ws.Range("Q16").Value2 = 408720
$ws.Range("R16").Value2 = 6703065
$ws.Range("S16").Value2 = 20
$ws.Range("AC16").Value2 = "Ringhack på gran"
$ws.Range("AW16").Value2 = "Moa Björnberg dillner"
$ws.Range("AX16").Value2 = "Moa Björnberg dillner"

# --- Row 17 becomes what Row 16 used to be ---
$ws.Range("A17").Value2 = 131063926
$ws.Range("B17").Value2 = 83090
$ws.Range("E17").Value2 = 1312
$ws.Range("F17").Value2 = "Gammelgransskål"
$ws.Range("G17").Value2 = "Pseudographis pinicola"
$ws.Range("H17").Value2 = "(Nyl.) Rehm"
$ws.Range("M17").Value2 = ""
$ws.Range("P17").Value2 = "Torsby kommun, Vrm"
$ws.Range("Q17").Value2 = 408603
$ws.Range("R17").Value2 = 6702927
$ws.Range("S17").Value2 = 5
$ws.Range("AC17").Value2 = ""
$ws.Range("AW17").Value2 = "Joakim Karlsson"
$ws.Range("AX17").Value2 = "Joakim Karlsson"
